# Adam Running Report - data refresh
# Organized some measures, Added more info to the read me file, and grouped some visuals.

$wb = $excel.ActiveWorkbook

$wsRuns  = $wb.Worksheets.Item("Runs")
$wsMiles = $wb.Worksheets.Item("Miles")
$wsShoes = $wb.Worksheets.Item("Shoes")
$wsGoals = $wb.Worksheets.Item("Goals")

# ---------------------------------------------------------------------------
# 1. Runs sheet: append two new run records (rows 43 and 44)
# ---------------------------------------------------------------------------
$wsRuns.Cells.Item(43, 1).Value = 42
$wsRuns.Cells.Item(43, 2).Value = 45018
$wsRuns.Cells.Item(43, 3).Value = 0.49652777777777773
$wsRuns.Cells.Item(43, 4).Value = "Sunday"
$wsRuns.Cells.Item(43, 5).Value = 4.8
$wsRuns.Cells.Item(43, 6).Value = 60.35
$wsRuns.Cells.Item(43, 7).Value = 12.26
$wsRuns.Cells.Item(43, 8).Value = 9.31
$wsRuns.Cells.Item(43, 9).Value = 645
$wsRuns.Cells.Item(43, 10).Value = 338
$wsRuns.Cells.Item(43, 12).Value = 5618
$wsRuns.Cells.Item(43, 13).Value = 163
$wsRuns.Cells.Item(43, 14).Value = 181
$wsRuns.Cells.Item(43, 15).Value = 137
$wsRuns.Cells.Item(43, 16).Value = 1
$wsRuns.Cells.Item(43, 17).Value = "Road"
$wsRuns.Cells.Item(43, 18).Value = 61
$wsRuns.Cells.Item(43, 19).Value = 39
$wsRuns.Cells.Item(43, 20).Value = "Colorado"
$wsRuns.Cells.Item(43, 21).Value = "Morrison"

$wsRuns.Cells.Item(44, 1).Value = 43
$wsRuns.Cells.Item(44, 2).Value = 45019
$wsRuns.Cells.Item(44, 3).Value = 0.40416666666666662
$wsRuns.Cells.Item(44, 4).Value = "Monday"
$wsRuns.Cells.Item(44, 5).Value = 4.01
$wsRuns.Cells.Item(44, 6).Value = 45.46
$wsRuns.Cells.Item(44, 7).Value = 11.25
$wsRuns.Cells.Item(44, 8).Value = 11.25
$wsRuns.Cells.Item(44, 9).Value = 587
$wsRuns.Cells.Item(44, 10).Value = 13
$wsRuns.Cells.Item(44, 12).Value = 5528
$wsRuns.Cells.Item(44, 13).Value = 160
$wsRuns.Cells.Item(44, 14).Value = 181
$wsRuns.Cells.Item(44, 15).Value = 146
$wsRuns.Cells.Item(44, 16).Value = 1
$wsRuns.Cells.Item(44, 17).Value = "Road"
$wsRuns.Cells.Item(44, 18).Value = 41
$wsRuns.Cells.Item(44, 19).Value = 39
$wsRuns.Cells.Item(44, 20).Value = "Colorado"
$wsRuns.Cells.Item(44, 21).Value = "Littleton"

# ---------------------------------------------------------------------------
# 2. Miles sheet: append per-mile splits for the two new runs (rows 191-200)
# ---------------------------------------------------------------------------
$milesRows = @(
    @(42, 1, 1,    10.33, 160),
    @(42, 2, 1,    13.18, 165),
    @(42, 3, 1,    13.19, 160),
    @(42, 4, 1,    12.18, 164),
    @(42, 5, 0.8,  13.43, 163),
    @(43, 1, 1,    10.47, 155),
    @(43, 2, 1,    12.19, 154),
    @(43, 3, 1,    11.13, 165),
    @(43, 4, 1,    11.18, 166),
    @(43, 5, 0.01, 9.44,  181)
)

$r = 191
foreach ($row in $milesRows) {
    $wsMiles.Cells.Item($r, 1).Value = $row[0]
    $wsMiles.Cells.Item($r, 2).Value = $row[1]
    $wsMiles.Cells.Item($r, 3).Value = $row[2]
    $wsMiles.Cells.Item($r, 4).Value = $row[3]
    $wsMiles.Cells.Item($r, 6).Value = $row[4]
    $r++
}

# ---------------------------------------------------------------------------
# 3. View state: freeze the header row on Runs and Miles, scroll each sheet
#    down to where the new data lives, and make Miles the active tab.
# ---------------------------------------------------------------------------
$wsRuns.Activate() | Out-Null
$wsRuns.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
$wsRuns.Rows.Item(1).EntireRow.Select() | Out-Null

$wsMiles.Activate() | Out-Null
$wsMiles.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
$wsMiles.Range("J193").Select() | Out-Null

$wsShoes.Activate() | Out-Null
$wsShoes.Range("C9").Select() | Out-Null

$wsGoals.Activate() | Out-Null
$wsGoals.Range("Q17").Select() | Out-Null

$wsMiles.Activate() | Out-Null

Write-Host "Workbook updated: added 2 runs, 10 mile splits, refreshed view state."
